$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.134.75"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.306.86"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.78"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.15"
$ws.Range("E6").Value = "  -3.02%  "
$ws.Range("E7").Value = "  +3.59%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.03"
$ws.Range("E13").Value = "  -3.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.88"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "2.665.41"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "2.217.11"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("D18").Value = "43.042.05"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.33"
$ws.Range("E19").Value = "  +7.78%  "
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.63"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.05"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.20"
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.81"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.18"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("E30").Value = "  -7.18%  "
$ws.Range("E31").Value = "  -6.18%  "
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.83"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.09"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.75"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("D42").Value = "2.004.92"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("E43").Value = "  -0.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.18"
$ws.Range("E44").Value = "  -3.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.16"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.30"
$ws.Range("E46").Value = "  -2.62%  "
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.38"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("D49").Value = "2.528.89"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.79"
$ws.Range("E51").Value = "  +11.38%  "
